# "Changes of 21st June 2022"
# Refresh a batch of FedEx tracking numbers (ShipmentTrackNum / PackageTrackNum)
# on Sheet1 rows 2-22, plus flip the Reject result in Q3 from PASS to FAIL.
#
# Note: a plain `$ws.Range($a1).Value = "320018538422"` would make Excel
# auto-coerce the long all-digit string into a *number*, which would change
# the cell's stored type (and silently drop any leading zeros/precision).
# Going through a quoted-literal formula and then "Copy / Paste Special -
# Values" collapses it back to a plain (non-formula) cell while keeping it
# text - exactly like the original workbook had it - without stamping a new
# cell style (no NumberFormat / quote-prefix detour needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$a1, [string]$text) {
    $cell = $ws.Range($a1)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue "C2"  "320018538422"
Set-TextValue "C3"  "320018538433"
Set-TextValue "Q3"  "FAIL"
Set-TextValue "C4"  "320018538466"
Set-TextValue "C5"  "320018538488"
Set-TextValue "D5"  "320018538488"
Set-TextValue "C6"  "320018538525"
Set-TextValue "D6"  "320018538525"
Set-TextValue "C7"  "320018538547"
Set-TextValue "D7"  "320018538547"
Set-TextValue "C8"  "320018538570"
Set-TextValue "C9"  "320018538591"
Set-TextValue "C10" "320018538628"
Set-TextValue "C11" "320018538640"
Set-TextValue "C12" "320018538694"
Set-TextValue "C13" "320018538710"
Set-TextValue "D13" "320018538710"
Set-TextValue "C14" "320018538742"
Set-TextValue "D14" "320018538742"
Set-TextValue "C15" "320018538775"
Set-TextValue "D15" "320018538775"
Set-TextValue "C16" "320018538801"
Set-TextValue "D16" "320018538801"
Set-TextValue "C17" "320018538823"
Set-TextValue "D17" "320018538823"
Set-TextValue "C18" "320018538867"
Set-TextValue "C19" "320018538889"
Set-TextValue "C20" "320018538915"
Set-TextValue "C21" "320018538937"
Set-TextValue "C22" "320018538960"
